# Auto-generated edit script: update market-price derived columns (H-N)
# on several leve-profit sheets, reflecting a scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1483.9286
$ws.Range("I32").Value = 816
$ws.Range("K32").Value = 816
$ws.Range("M32").Value = -490

$ws = $wb.Worksheets.Item("ALC")
# Row 106
$ws.Range("H106").Value = 2743.1333
$ws.Range("I106").Value = 2797.2144
$ws.Range("K106").Value = 2797.2144
$ws.Range("M106").Value = -2166.2144

$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 3219.4
$ws.Range("I113").Value = 2400
$ws.Range("J113").Value = 3424.25
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 3424.25
$ws.Range("M113").Value = 854
$ws.Range("N113").Value = -9932.25

$ws = $wb.Worksheets.Item("ALC")
# Row 138
$ws.Range("H138").Value = 495977.84
$ws.Range("I138").Value = 999.8182
$ws.Range("J138").Value = 566689
$ws.Range("K138").Value = 2999.4546
$ws.Range("L138").Value = 1700067
$ws.Range("M138").Value = 2140.5454
$ws.Range("N138").Value = -1710347

$ws = $wb.Worksheets.Item("ALC")
# Row 141
$ws.Range("H141").Value = 715.7143
$ws.Range("I141").Value = 715.7143
$ws.Range("K141").Value = 2147.1429
$ws.Range("M141").Value = 3032.8571

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 6781.1177
$ws.Range("I2").Value = 929.5
$ws.Range("K2").Value = 929.5
$ws.Range("M2").Value = -816.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2015.1111
$ws.Range("J61").Value = 2833.3333
$ws.Range("L61").Value = 2833.3333
$ws.Range("N61").Value = -3257.3333

$ws = $wb.Worksheets.Item("ARM")
# Row 116
$ws.Range("H116").Value = 6781.1177
$ws.Range("I116").Value = 929.5
$ws.Range("K116").Value = 929.5
$ws.Range("M116").Value = 1364.5

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 1638.6666
$ws.Range("I122").Value = 1666.3334
$ws.Range("K122").Value = 4999.0002
$ws.Range("M122").Value = -2549.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 136
$ws.Range("H136").Value = 2015.1111
$ws.Range("J136").Value = 2833.3333
$ws.Range("L136").Value = 8499.999899999999
$ws.Range("N136").Value = -13599.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 6781.1177
$ws.Range("I3").Value = 929.5
$ws.Range("K3").Value = 929.5
$ws.Range("M3").Value = -815.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 16667513
$ws.Range("I94").Value = 20834010
$ws.Range("K94").Value = 20834010
$ws.Range("M94").Value = -20833559

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1567.1111
$ws.Range("I107").Value = 1157
$ws.Range("J107").Value = 2264.3
$ws.Range("K107").Value = 1157
$ws.Range("L107").Value = 2264.3
$ws.Range("M107").Value = 763
$ws.Range("N107").Value = -6104.3

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 4541.483
$ws.Range("I134").Value = 1100.1482
$ws.Range("J134").Value = 50999.5
$ws.Range("K134").Value = 3300.4446
$ws.Range("L134").Value = 152998.5
$ws.Range("M134").Value = -765.4446000000003
$ws.Range("N134").Value = -158068.5

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 10400
$ws.Range("I23").Value = 4000
$ws.Range("K23").Value = 4000
$ws.Range("M23").Value = -3760

$ws = $wb.Worksheets.Item("CRP")
# Row 27
$ws.Range("H27").Value = 10400
$ws.Range("I27").Value = 4000
$ws.Range("K27").Value = 4000
$ws.Range("M27").Value = -3808

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1599.7142
$ws.Range("I31").Value = 1616.3334
$ws.Range("K31").Value = 1616.3334
$ws.Range("M31").Value = -1321.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 34
$ws.Range("H34").Value = 1599.7142
$ws.Range("I34").Value = 1616.3334
$ws.Range("K34").Value = 1616.3334
$ws.Range("M34").Value = -1414.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 1722.6666
$ws.Range("I99").Value = 1667.2
$ws.Range("K99").Value = 1667.2
$ws.Range("M99").Value = -169.2

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 883
$ws.Range("I105").Value = 883
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 883
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 864
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 657.2273
$ws.Range("I107").Value = 432.1
$ws.Range("J107").Value = 844.8333
$ws.Range("K107").Value = 432.1
$ws.Range("L107").Value = 844.8333
$ws.Range("M107").Value = 1487.9
$ws.Range("N107").Value = -4684.8333

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 1047.7778
$ws.Range("I122").Value = 1004.2857
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3012.8571
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -562.8571000000002
$ws.Range("N122").Value = -8500

$ws = $wb.Worksheets.Item("CRP")
# Row 126
$ws.Range("H126").Value = 1722.6666
$ws.Range("I126").Value = 1667.2
$ws.Range("K126").Value = 5001.6
$ws.Range("M126").Value = -2531.6

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 5846.3335
$ws.Range("I132").Value = 6295.381
$ws.Range("J132").Value = 4274.6665
$ws.Range("K132").Value = 18886.143
$ws.Range("L132").Value = 12823.9995
$ws.Range("M132").Value = -16356.143
$ws.Range("N132").Value = -17883.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 185.1
$ws.Range("I38").Value = 93.875
$ws.Range("K38").Value = 281.625
$ws.Range("M38").Value = 65.375

$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 15000.111
$ws.Range("J63").Value = 15000.111
$ws.Range("L63").Value = 15000.111
$ws.Range("N63").Value = -16372.111

$ws = $wb.Worksheets.Item("GSM")
# Row 66
$ws.Range("H66").Value = 15000.111
$ws.Range("J66").Value = 15000.111
$ws.Range("L66").Value = 45000.333
$ws.Range("N66").Value = -51864.333

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1465
$ws.Range("I122").Value = 1430
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4290
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1840
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 35718084
$ws.Range("I122").Value = 62503396
$ws.Range("K122").Value = 187510188
$ws.Range("M122").Value = -187507738

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 34037.516
$ws.Range("I132").Value = 1484.5
$ws.Range("K132").Value = 4453.5
$ws.Range("M132").Value = -1923.5

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 546.63635
$ws.Range("I107").Value = 517.8333
$ws.Range("J107").Value = 581.2
$ws.Range("K107").Value = 1553.4999
$ws.Range("L107").Value = 1743.6
$ws.Range("M107").Value = 366.5001
$ws.Range("N107").Value = -5583.6

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 20006902
$ws.Range("I122").Value = 21673894
$ws.Range("K122").Value = 65021682
$ws.Range("M122").Value = -65019232

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3135.6924
$ws.Range("I132").Value = 2501.6316
$ws.Range("J132").Value = 4856.7144
$ws.Range("K132").Value = 7504.8948
$ws.Range("L132").Value = 14570.1432
$ws.Range("M132").Value = -4974.8948
$ws.Range("N132").Value = -19630.1432

$ws = $wb.Worksheets.Item("WVR")
# Row 133
$ws.Range("H133").Value = 35531.668
$ws.Range("J133").Value = 35531.668
$ws.Range("L133").Value = 35531.668
$ws.Range("N133").Value = -45651.668
